$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ------------------------------------------------------------------
# Sheet "Templates": insert a new row for the "Apology" hotstring,
# refresh the other rows' text, and append a "Test NoPrefix Hotkey"
# row at the end.
# ------------------------------------------------------------------
$ws1.Activate()

# Insert a new row at 3 - shifts old rows 3,4,5 down to 4,5,6.
$ws1.Rows.Item(3).Insert()

# Row 2 id stays 0 (text refreshed further down)
$ws1.Range("A2").Value = 0

# Row 3 - Apology (new hotstring)
$ws1.Range("A3").Value = 1
$ws1.Range("C3").Value = "<sorry"
$ws1.Range("D3").Value = "I am sorry for the inconvenience."

# Row 4 - Greeting, now using the built-in <hi prefix
$ws1.Range("A4").Value = 2
$ws1.Range("C4").Value = "<hi"
$ws1.Range("D4").Value = "Hello World"

# Row 6 - Test NoPrefix Hotkey (new row, appended after Grocery Note)
$ws1.Range("A6").Value = 4
$ws1.Range("C6").Value = "test"
$ws1.Range("D6").Value = "This is a hotstring without a prefix."
$ws1.Range("B6").Value = "Test NoPrefix Hotkey"

$ws1.Range("B3").Value = "Apology"

# Row 5 - Grocery Note (content unchanged, just shifted down by the insert)
$ws1.Range("A5").Value = 3
$ws1.Range("B5").Value = "Grocery Note"
$ws1.Range("C5").Value = "<list"
$ws1.Range("D5").Value = "This is the list of items I need from the store. `n`n    * Apples`n    * Oranges`n    * Paper Towels`n`nThat’s it, the end of the list."

# Row 2 - Introduction (renamed last)
$ws1.Range("B2").Value = "Introduction"
$ws1.Range("C2").Value = "<ate"
$ws1.Range("D2").Value = "AutoHotkey Text Expander"

# Drop the leftover column/cell formatting inherited from the old layout.
$ws1.Cells.ClearFormats()

# Re-apply word-wrap to the "Extended Text" cells that need it; the long
# grocery note re-triggers automatic row auto-fit (~100.8pt).
$ws1.Range("D4").WrapText = $true
$ws1.Range("D5").WrapText = $true
$ws1.Range("D9").WrapText = $true

$ws1.Range("A1").Select()

# ------------------------------------------------------------------
# Sheet "Info": move the active selection down to A2.
# ------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A2").Select()

# "Templates" is the tab that should end up selected/active.
$ws1.Activate()
